# "add neutrons final pictures"
#
# Лист1 gets one more measured point: V_0 = 13.7 kV (new row 5), with the
# same formulas (copied down from row 4) computing U_THGEM_0, Epsilon_GAr,
# Epsilon_LAr_emiss and Epsilon_LAr_drift for that point.
#
# The workbook's active/selected sheet also moves from Лист2 back to Лист1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")
$ws2 = $wb.Worksheets.Item("Лист2")

# --- Лист1: add row 5 with the new measurement point (A5 = 13.7) -------
$ws1.Range("A5").Value = 13.7

$ws1.Range("C5").Formula = "=A5/`$K`$2*600"
$ws1.Range("D5").Formula = "=C5 / (1.8 + 0.4 / `$J`$2)"
$ws1.Range("E5").Formula = "=D5/`$J`$2"
$ws1.Range("F5").Formula = "=(A5/`$K`$2*40) / 1.6"

# reuse the same number format ("0.00", style index 1) as the rows above
$ws1.Range("C5:F5").NumberFormat = $ws1.Range("C4:F4").NumberFormat

# --- move the selected/active sheet from Лист2 back to Лист1 -----------
# update Лист2's own selection first, while it is still the active sheet
$ws2.Range("E3").Select()

# activate Лист1 last so that it ends up being the workbook's active tab
$ws1.Activate()
$ws1.Range("E9").Select()
